# Update cryptos list with latest price/volume data
# (row 17/18 and 40/41 pairs swap rank order; row 51 swaps EOS for BabyDogeCoin)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.846.83'
$ws.Range('E2').Value = '  +1.74%  '
$ws.Range('D3').Value = '''1.886.10'
$ws.Range('E3').Value = '  +1.45%  '
$ws.Range('D4').Value = '''1.010'
$ws.Range('E4').Value = '  +0.88%  '
$ws.Range('D5').Value = '''333.91'
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('D6').Value = '''1.008'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('D7').Value = '''0.4741'
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = '''0.3932'
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('D9').Value = '''47.73'
$ws.Range('E9').Value = '  +1.22%  '
$ws.Range('D10').Value = '''0.08061'
$ws.Range('E10').Value = '  +0.72%  '
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('D12').Value = '''22.12'
$ws.Range('E12').Value = '  +2.61%  '
$ws.Range('D13').Value = '''1.890.16'
$ws.Range('E13').Value = '  +2.02%  '
$ws.Range('D14').Value = '''5.986'
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').Value = '''7.150'
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('D16').Value = '''1.011'
$ws.Range('E16').Value = '  +0.74%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = '''0.06736'
$ws.Range('E17').Value = '  +2.92%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '''0.00001052'
$ws.Range('E18').Value = '  +1.46%  '
$ws.Range('D19').Value = '''87.43'
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('D20').Value = '''17.24'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('D22').Value = '''27.872.10'
$ws.Range('E22').Value = '  +1.84%  '
$ws.Range('D23').Value = '''5.520'
$ws.Range('E23').Value = '  +0.43%  '
$ws.Range('D24').Value = '''10.98'
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('D25').Value = '''2.333'
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('D26').Value = '''2.110.37'
$ws.Range('E26').Value = '  +1.74%  '
$ws.Range('D27').Value = '''159.40'
$ws.Range('E27').Value = '  +3.42%  '
$ws.Range('D28').Value = '''20.10'
$ws.Range('E28').Value = '  -1.48%  '
$ws.Range('D29').Value = '''2.106'
$ws.Range('E29').Value = '  +1.36%  '
$ws.Range('D30').Value = '''5.555'
$ws.Range('E30').Value = '  +0.87%  '
$ws.Range('D31').Value = '''121.94'
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('D32').Value = '''0.9784'
$ws.Range('E32').Value = '  +2.42%  '
$ws.Range('D33').Value = '''0.09500'
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('E34').Value = '  +0.24%  '
$ws.Range('D35').Value = '''3.638'
$ws.Range('E35').Value = '  +1.37%  '
$ws.Range('D36').Value = '''5.362'
$ws.Range('E36').Value = '  +1.42%  '
$ws.Range('E37').Value = '  +1.80%  '
$ws.Range('D38').Value = '''0.02273'
$ws.Range('E38').Value = '  +1.68%  '
$ws.Range('D39').Value = '''1.222'
$ws.Range('E39').Value = '  +1.06%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '''0.6012'
$ws.Range('E40').Value = '  +1.02%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '''8.067'
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('D42').Value = '''0.1899'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').Value = '''10.29'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').Value = '''1.260'
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('D45').Value = '''0.5709'
$ws.Range('E45').Value = '  +0.81%  '
$ws.Range('D46').Value = '''12.29'
$ws.Range('E46').Value = '  +1.56%  '
$ws.Range('D47').Value = '''3.403'
$ws.Range('E47').Value = '  -0.78%  '
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('D49').Value = '''0.06916'
$ws.Range('E49').Value = '  +2.13%  '
$ws.Range('D50').Value = '''113.57'
$ws.Range('E50').Value = '  +3.71%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '''0.00000000300'
$ws.Range('E51').Value = '  +6.02%  '
